# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the freshly generated output data.

$wb = $excel.ActiveWorkbook

# Row -> (old value, new value) updates for worksheet "展览" (sheet1)
$zhanlanUpdates = @{
    6  = 1062
    8  = 8080
    9  = 132
    10 = 199
    11 = 6861
    14 = 4938
    17 = 5352
    18 = 1070
    19 = 324
    21 = 448
    27 = 9073
    29 = 1631
    33 = 841
    37 = 1174
    39 = 4727
}

# Row -> new value updates for worksheet "全部类型" (sheet4)
$quanbuUpdates = @{
    8  = 1062
    10 = 8080
    11 = 132
    12 = 199
    13 = 6861
    17 = 4938
    19 = 5352
    20 = 1070
    21 = 324
    23 = 448
    30 = 9073
    32 = 1631
    35 = 841
    39 = 1174
    41 = 4727
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
